$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 formatting (row height) to match the rest of the question rows
$ws.Rows.Item(3).RowHeight = 15.5

# Add new question/answer row at row 19
$ws.Range("C19").Value = "who is most hardworker in team"
$ws.Range("C19").Font.Name = "Lato"
$ws.Range("C19").Font.Size = 10
$ws.Range("C19").Font.Color = 5188908
$ws.Range("D19").Value = "Max vote wins"
$ws.Rows.Item(19).RowHeight = 15.5

$ws.Range("C14").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
